$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'304.76"
$ws.Range("E2").Value = "'2.23%"
$ws.Range("D3").Value = "'35.66"
$ws.Range("E3").Value = "'12.72%"
$ws.Range("D4").Value = "'5.095"
$ws.Range("E4").Value = "'2.09%"
$ws.Range("D5").Value = "'0.07807"
$ws.Range("E5").Value = "'1.56%"
$ws.Range("D6").Value = "'2.265"
$ws.Range("E6").Value = "'0.83%"
$ws.Range("D7").Value = "'8.107"
$ws.Range("E7").Value = "'3.03%"
$ws.Range("D8").Value = "'4.038"
$ws.Range("E8").Value = "'6.40%"
$ws.Range("D9").Value = "'0.9283"
$ws.Range("E9").Value = "'0.42%"
$ws.Range("D10").Value = "'0.09749"
$ws.Range("E10").Value = "'-0.87%"
$ws.Range("D11").Value = "'0.1830"
$ws.Range("E11").Value = "'4.82%"
$ws.Range("D12").Value = "'0.08618"
$ws.Range("E12").Value = "'2.60%"
$ws.Range("D13").Value = "'0.03415"
$ws.Range("E13").Value = "'5.10%"
$ws.Range("D14").Value = "'0.09943"
$ws.Range("E14").Value = "'1.18%"
$ws.Range("D15").Value = "'0.001474"
$ws.Range("E15").Value = "'-0.41%"
$ws.Range("D16").Value = "'0.005674"
$ws.Range("E16").Value = "'-1.54%"
$ws.Range("D17").Value = "'3.480"
$ws.Range("E17").Value = "'-1.18%"
$ws.Range("E18").Value = "'-0.81%"
$ws.Range("E19").Value = "'2.90%"
$ws.Range("D20").Value = "'0.1324"
$ws.Range("E20").Value = "'-0.20%"
$ws.Range("E21").Value = "'11.76%"
$ws.Range("D22").Value = "'0.2241"
$ws.Range("E22").Value = "'-1.60%"
$ws.Range("D23").Value = "'0.04689"
$ws.Range("E23").Value = "'4.37%"
$ws.Range("E24").Value = "'2.54%"
$ws.Range("D25").Value = "'0.004539"
$ws.Range("E25").Value = "'4.11%"
$ws.Range("D26").Value = "'0.0001302"
$ws.Range("E26").Value = "'1.34%"
$ws.Range("E27").Value = "'-19.81%"
$ws.Range("D39").Value = "'0.01763"
$ws.Range("E39").Value = "'4.43%"
$ws.Range("D40").Value = "'0.04717"
$ws.Range("E40").Value = "'1.73%"
$ws.Range("D41").Value = "'0.007922"
$ws.Range("E41").Value = "'5.11%"
$ws.Range("D42").Value = "'0.1419"
$ws.Range("E42").Value = "'2.38%"
$ws.Range("D43").Value = "'0.008010"
$ws.Range("E43").Value = "'-17.70%"
$ws.Range("D44").Value = "'0.002294"
$ws.Range("E44").Value = "'11.17%"
$ws.Range("D45").Value = "'0.009106"
$ws.Range("E45").Value = "'-12.95%"
$ws.Range("D46").Value = "'0.00006214"
$ws.Range("E46").Value = "'2.48%"
$ws.Range("E47").Value = "'1.01%"
$ws.Range("D48").Value = "'4.040"
$ws.Range("E48").Value = "'58.39%"
$ws.Range("D49").Value = "'0.002696"
$ws.Range("E49").Value = "'35.86%"
$ws.Range("D50").Value = "'0.00002104"
$ws.Range("E50").Value = "'1.01%"
$ws.Range("D51").Value = "'0.0002004"
$ws.Range("E51").Value = "'1.01%"
